# Auto-generated Excel COM-interop script
# Applies bulk market-price value updates across multiple sheets
# (Leve profit calculations refreshed by scheduled runner)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 48354.5
$ws.Range("I64").Value = 69359.92999999999
$ws.Range("K64").Value = 69359.92999999999
$ws.Range("M64").Value = -69111.92999999999
$ws.Range("H67").Value = 48354.5
$ws.Range("I67").Value = 69359.92999999999
$ws.Range("K67").Value = 69359.92999999999
$ws.Range("M67").Value = -68501.92999999999
$ws.Range("H74").Value = 3028.7896
$ws.Range("I74").Value = 2824.8333
$ws.Range("J74").Value = 3378.4285
$ws.Range("K74").Value = 2824.8333
$ws.Range("L74").Value = 3378.4285
$ws.Range("M74").Value = -1888.8333
$ws.Range("N74").Value = -5250.4285
$ws.Range("H76").Value = 5626.125
$ws.Range("I76").Value = 5003
$ws.Range("J76").Value = 5715.143
$ws.Range("K76").Value = 5003
$ws.Range("L76").Value = 5715.143
$ws.Range("M76").Value = -4688
$ws.Range("N76").Value = -6345.143
$ws.Range("H77").Value = 3028.7896
$ws.Range("I77").Value = 2824.8333
$ws.Range("J77").Value = 3378.4285
$ws.Range("K77").Value = 14124.1665
$ws.Range("L77").Value = 16892.1425
$ws.Range("M77").Value = -9444.166499999999
$ws.Range("N77").Value = -26252.1425
$ws.Range("H79").Value = 5626.125
$ws.Range("I79").Value = 5003
$ws.Range("J79").Value = 5715.143
$ws.Range("K79").Value = 5003
$ws.Range("L79").Value = 5715.143
$ws.Range("M79").Value = -3911
$ws.Range("N79").Value = -7899.143
$ws.Range("H131").Value = 4392.2354
$ws.Range("I131").Value = 1113.4286
$ws.Range("J131").Value = 4913.864
$ws.Range("K131").Value = 3340.2858
$ws.Range("L131").Value = 14741.592
$ws.Range("M131").Value = 1699.7142
$ws.Range("N131").Value = -24821.592
$ws.Range("H138").Value = 3773.577
$ws.Range("I138").Value = 3873.75
$ws.Range("J138").Value = 3762.1287
$ws.Range("K138").Value = 11621.25
$ws.Range("L138").Value = 11286.3861
$ws.Range("M138").Value = -6481.25
$ws.Range("N138").Value = -21566.3861

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 168397.67
$ws.Range("I105").Value = 168463.33
$ws.Range("J105").Value = 168332
$ws.Range("K105").Value = 168463.33
$ws.Range("L105").Value = 168332
$ws.Range("M105").Value = -166716.33
$ws.Range("N105").Value = -171826

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1185.3334
$ws.Range("I16").Value = 1174.75
$ws.Range("K16").Value = 1174.75
$ws.Range("M16").Value = -887.75
$ws.Range("H31").Value = 48034.312
$ws.Range("I31").Value = 1534.9333
$ws.Range("J31").Value = 89063.17999999999
$ws.Range("K31").Value = 1534.9333
$ws.Range("L31").Value = 89063.17999999999
$ws.Range("M31").Value = -1239.9333
$ws.Range("N31").Value = -89653.17999999999
$ws.Range("H34").Value = 48034.312
$ws.Range("I34").Value = 1534.9333
$ws.Range("J34").Value = 89063.17999999999
$ws.Range("K34").Value = 1534.9333
$ws.Range("L34").Value = 89063.17999999999
$ws.Range("M34").Value = -1332.9333
$ws.Range("N34").Value = -89467.17999999999
$ws.Range("H62").Value = 2371.0417
$ws.Range("I62").Value = 2243.75
$ws.Range("J62").Value = 2625.625
$ws.Range("K62").Value = 2243.75
$ws.Range("L62").Value = 2625.625
$ws.Range("M62").Value = -1619.75
$ws.Range("N62").Value = -3873.625
$ws.Range("H65").Value = 2371.0417
$ws.Range("I65").Value = 2243.75
$ws.Range("J65").Value = 2625.625
$ws.Range("K65").Value = 11218.75
$ws.Range("L65").Value = 13128.125
$ws.Range("M65").Value = -8098.75
$ws.Range("N65").Value = -19368.125
$ws.Range("H113").Value = 1185.3334
$ws.Range("I113").Value = 1174.75
$ws.Range("K113").Value = 1174.75
$ws.Range("M113").Value = 995.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1243.9678
$ws.Range("I5").Value = 878.03705
$ws.Range("J5").Value = 1526.2572
$ws.Range("K5").Value = 2634.11115
$ws.Range("L5").Value = 4578.7716
$ws.Range("M5").Value = -2522.11115
$ws.Range("N5").Value = -4802.7716
$ws.Range("H134").Value = 4453.52
$ws.Range("I134").Value = 2463.2
$ws.Range("J134").Value = 7439
$ws.Range("K134").Value = 7389.599999999999
$ws.Range("L134").Value = 22317
$ws.Range("M134").Value = -2319.599999999999
$ws.Range("N134").Value = -32457
$ws.Range("H135").Value = 1243.9678
$ws.Range("I135").Value = 878.03705
$ws.Range("J135").Value = 1526.2572
$ws.Range("K135").Value = 7902.33345
$ws.Range("L135").Value = 13736.3148
$ws.Range("M135").Value = -5367.33345
$ws.Range("N135").Value = -18806.3148

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2038.3334
$ws.Range("J43").Value = 5865.2856
$ws.Range("L43").Value = 5865.2856
$ws.Range("N43").Value = -6167.2856
$ws.Range("H63").Value = 14600
$ws.Range("J63").Value = 17800
$ws.Range("L63").Value = 17800
$ws.Range("N63").Value = -19172
$ws.Range("H66").Value = 14600
$ws.Range("J66").Value = 17800
$ws.Range("L66").Value = 53400
$ws.Range("N66").Value = -60264
$ws.Range("H70").Value = 130951.375
$ws.Range("I70").Value = 159401.69
$ws.Range("K70").Value = 159401.69
$ws.Range("M70").Value = -159131.69
$ws.Range("H73").Value = 130951.375
$ws.Range("I73").Value = 159401.69
$ws.Range("K73").Value = 159401.69
$ws.Range("M73").Value = -158465.69
$ws.Range("H80").Value = 167037920
$ws.Range("I80").Value = 250551250
$ws.Range("J80").Value = 11250
$ws.Range("K80").Value = 250551250
$ws.Range("L80").Value = 11250
$ws.Range("M80").Value = -250550252
$ws.Range("N80").Value = -13246
$ws.Range("H83").Value = 167037920
$ws.Range("I83").Value = 250551250
$ws.Range("J83").Value = 11250
$ws.Range("K83").Value = 1252756250
$ws.Range("L83").Value = 56250
$ws.Range("M83").Value = -1252751258
$ws.Range("N83").Value = -66234

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1209.2142
$ws.Range("I55").Value = 1329.3077
$ws.Range("J55").Value = 1155.3793
$ws.Range("K55").Value = 1329.3077
$ws.Range("L55").Value = 1155.3793
$ws.Range("M55").Value = -1156.3077
$ws.Range("N55").Value = -1501.3793

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716
$ws.Range("H107").Value = 111699.555
$ws.Range("I107").Value = 650
$ws.Range("J107").Value = 200539.2
$ws.Range("K107").Value = 1950
$ws.Range("L107").Value = 601617.6000000001
$ws.Range("M107").Value = -30
$ws.Range("N107").Value = -605457.6000000001

Write-Host "Applied all leve profit value updates."